$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (and its entry in workbook.xml)
$ws.Name = "Tab_3a_Postulate"

# New rows 41:43 need to be created below the former last row (40).
# Copy formatting from an existing data row down into the new rows first,
# so the new cells pick up the same style used by the rest of the table.
$ws.Range("A2:D2").Copy()
$ws.Range("A41:D43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update header row: add a PNr column, shift BNr into column A, drop ZNr
$ws.Cells.Item(1,1).Value = "PNr"
$ws.Cells.Item(1,2).Value = "BNr"

# Update column widths to the new (wider) layout
$ws.Columns.Item(1).ColumnWidth = 13.072
$ws.Columns.Item(2).ColumnWidth = 10.929
$ws.Columns.Item(3).ColumnWidth = 66.786
$ws.Columns.Item(4).ColumnWidth = 70.072

# Populate data rows 2..43 with the new Postulate-level data (PNr, BNr, BezDe, BezEn)
$ws.Cells.Item(2,1).Value = "Z01_B01_P01"
$ws.Cells.Item(2,2).Value = "Z01_B01"
$ws.Cells.Item(2,3).Value = "Armut begrenzen"
$ws.Cells.Item(2,4).Value = "Limiting poverty"

$ws.Cells.Item(3,1).Value = "Z02_B01_P01"
$ws.Cells.Item(3,2).Value = "Z02_B01"
$ws.Cells.Item(3,3).Value = "In unseren Kulturlandschaften umweltverträglich produzieren"
$ws.Cells.Item(3,4).Value = "Environmentally sound production in our cultivated landscapes"

$ws.Cells.Item(4,1).Value = "Z02_B02_P01"
$ws.Cells.Item(4,2).Value = "Z02_B02"
$ws.Cells.Item(4,3).Value = "Das Recht auf Nahrung weltweit verwirklichen"
$ws.Cells.Item(4,4).Value = "Realising globally the right to food"

$ws.Cells.Item(5,1).Value = "Z03_B01_P01"
$ws.Cells.Item(5,2).Value = "Z03_B01"
$ws.Cells.Item(5,3).Value = "Länger gesund leben"
$ws.Cells.Item(5,4).Value = "Living healthy longer"

$ws.Cells.Item(6,1).Value = "Z03_B02_P01"
$ws.Cells.Item(6,2).Value = "Z03_B02"
$ws.Cells.Item(6,3).Value = "Gesunde Umwelt erhalten"
$ws.Cells.Item(6,4).Value = "Keeping the environment healthy"

$ws.Cells.Item(7,1).Value = "Z03_B03_P01"
$ws.Cells.Item(7,2).Value = "Z03_B03"
$ws.Cells.Item(7,3).Value = "Globale Gesundheitsarchitektur stärken"
$ws.Cells.Item(7,4).Value = "Strengthening the global health architecture"

$ws.Cells.Item(8,1).Value = "Z04_B01_P01"
$ws.Cells.Item(8,2).Value = "Z04_B01"
$ws.Cells.Item(8,3).Value = "Bildung und Qualifikation kontinuierlich verbessern"
$ws.Cells.Item(8,4).Value = "Continuously improving education and vocational training"

$ws.Cells.Item(9,1).Value = "Z04_B02_P01"
$ws.Cells.Item(9,2).Value = "Z04_B02"
$ws.Cells.Item(9,3).Value = "Vereinbarkeit von Familie und Beruf verbessern"
$ws.Cells.Item(9,4).Value = "Improving the compatibility of work and family life"

$ws.Cells.Item(10,1).Value = "Z05_B01_P01"
$ws.Cells.Item(10,2).Value = "Z05_B01"
$ws.Cells.Item(10,3).Value = "Gleichstellung und partnerschaftliche Aufgabenteilung fördern"
$ws.Cells.Item(10,4).Value = "Promoting equal opportunities in society"

$ws.Cells.Item(11,1).Value = "Z05_B01_P02"
$ws.Cells.Item(11,2).Value = "Z05_B01"
$ws.Cells.Item(11,3).Value = "Wirtschaftliche Teilhabe von Frauen global stärken"
$ws.Cells.Item(11,4).Value = "Strengthening the economic participation of women globally"

$ws.Cells.Item(12,1).Value = "Z06_B01_P01"
$ws.Cells.Item(12,2).Value = "Z06_B01"
$ws.Cells.Item(12,3).Value = "Minderung der stofflichen Belastung von Gewässern"
$ws.Cells.Item(12,4).Value = "Reducing the pollution of water with substances"

$ws.Cells.Item(13,1).Value = "Z06_B02_P01"
$ws.Cells.Item(13,2).Value = "Z06_B02"
$ws.Cells.Item(13,3).Value = "Besserer Zugang zu Trinkwasser- und Sanitärversorgung weltweit, höhere (sichere) Qualität"
$ws.Cells.Item(13,4).Value = "Better access to drinking water and sanitation worldwide, higher (safer) quality"

$ws.Cells.Item(14,1).Value = "Z07_B01_P01"
$ws.Cells.Item(14,2).Value = "Z07_B01"
$ws.Cells.Item(14,3).Value = "Ressourcen sparsam und effizient nutzen"
$ws.Cells.Item(14,4).Value = "Using resources economically and efficiently"

$ws.Cells.Item(15,1).Value = "Z07_B02_P01"
$ws.Cells.Item(15,2).Value = "Z07_B02"
$ws.Cells.Item(15,3).Value = "Zukunftsfähige Energieversorgung ausbauen"
$ws.Cells.Item(15,4).Value = "Strengthening a sustainable energy supply"

$ws.Cells.Item(16,1).Value = "Z08_B01_P01"
$ws.Cells.Item(16,2).Value = "Z08_B01"
$ws.Cells.Item(16,3).Value = "Ressourcen sparsam und effizient nutzen"
$ws.Cells.Item(16,4).Value = "Using resources economically and efficiently"

$ws.Cells.Item(17,1).Value = "Z08_B02_P01"
$ws.Cells.Item(17,2).Value = "Z08_B02"
$ws.Cells.Item(17,3).Value = "Staatsfinanzen konsolidieren – Generationengerechtigkeit schaffen"
$ws.Cells.Item(17,4).Value = "Consolidating public finances – Creating intergenerational equity"

$ws.Cells.Item(18,1).Value = "Z08_B03_P01"
$ws.Cells.Item(18,2).Value = "Z08_B03"
$ws.Cells.Item(18,3).Value = "Gute Investitionsbedingungen schaffen – Wohlstand dauerhaft erhalten"
$ws.Cells.Item(18,4).Value = "Creating favourable investment conditions – Securing long-term prosperity"

$ws.Cells.Item(19,1).Value = "Z08_B04_P01"
$ws.Cells.Item(19,2).Value = "Z08_B04"
$ws.Cells.Item(19,3).Value = "Wirtschaftsleistung umwelt- und sozialverträglich steigern"
$ws.Cells.Item(19,4).Value = "Combining greater economic output with environmental and social responsibility"

$ws.Cells.Item(20,1).Value = "Z08_B05_P01"
$ws.Cells.Item(20,2).Value = "Z08_B05"
$ws.Cells.Item(20,3).Value = "Beschäftigungsniveau steigern"
$ws.Cells.Item(20,4).Value = "Boosting employment levels"

$ws.Cells.Item(21,1).Value = "Z08_B06_P01"
$ws.Cells.Item(21,2).Value = "Z08_B06"
$ws.Cells.Item(21,3).Value = "Menschenwürdige Arbeit weltweit ermöglichen"
$ws.Cells.Item(21,4).Value = "Enabling decent work worldwide"

$ws.Cells.Item(22,1).Value = "Z09_B01_P01"
$ws.Cells.Item(22,2).Value = "Z09_B01"
$ws.Cells.Item(22,3).Value = "Zukunft mit neuen Lösungen nachhaltig gestalten"
$ws.Cells.Item(22,4).Value = "Shaping the future with new solutions"

$ws.Cells.Item(23,1).Value = "Z10_B01_P01"
$ws.Cells.Item(23,2).Value = "Z10_B01"
$ws.Cells.Item(23,3).Value = "Schulische Bildungserfolge von Ausländern in Deutschland verbessern"
$ws.Cells.Item(23,4).Value = "Improving educational success of foreigners in German schools"

$ws.Cells.Item(24,1).Value = "Z10_B02_P01"
$ws.Cells.Item(24,2).Value = "Z10_B02"
$ws.Cells.Item(24,3).Value = "Zu große Ungleichheit innerhalb Deutschlands verhindern"
$ws.Cells.Item(24,4).Value = "Preventing excessive inequality within Germany"

$ws.Cells.Item(25,1).Value = "Z11_B01_P01"
$ws.Cells.Item(25,2).Value = "Z11_B01"
$ws.Cells.Item(25,3).Value = "Flächen nachhaltig nutzen"
$ws.Cells.Item(25,4).Value = "Using land sustainably"

$ws.Cells.Item(26,1).Value = "Z11_B02_P01"
$ws.Cells.Item(26,2).Value = "Z11_B02"
$ws.Cells.Item(26,3).Value = "Mobilität sichern – Umwelt schonen"
$ws.Cells.Item(26,4).Value = "Guaranteeing mobility – Protecting the environment"

$ws.Cells.Item(27,1).Value = "Z11_B03_P01"
$ws.Cells.Item(27,2).Value = "Z11_B03"
$ws.Cells.Item(27,3).Value = "Bezahlbarer Wohnraum für alle"
$ws.Cells.Item(27,4).Value = "Affordable housing for all"

$ws.Cells.Item(28,1).Value = "Z11_B04_P01"
$ws.Cells.Item(28,2).Value = "Z11_B04"
$ws.Cells.Item(28,3).Value = "Zugang zum Kulturerbe verbessern"
$ws.Cells.Item(28,4).Value = "Improving access to cultural heritage"

$ws.Cells.Item(29,1).Value = "Z12_B01_P01"
$ws.Cells.Item(29,2).Value = "Z12_B01"
$ws.Cells.Item(29,3).Value = "Konsum umwelt- und sozialverträglich gestalten"
$ws.Cells.Item(29,4).Value = "Making consumption environmentally and socially compatible"

$ws.Cells.Item(30,1).Value = "Z12_B02_P01"
$ws.Cells.Item(30,2).Value = "Z12_B02"
$ws.Cells.Item(30,3).Value = "Anteil nachhaltiger Produktion stetig erhöhen"
$ws.Cells.Item(30,4).Value = "Increasing the proportion of sustainable production continuously"

$ws.Cells.Item(31,1).Value = "Z12_B03_P01"
$ws.Cells.Item(31,2).Value = "Z12_B03"
$ws.Cells.Item(31,3).Value = "Vorbildwirkung der öffentlichen Hand für nachhaltige öffentliche Beschaffung verwirklichen"
$ws.Cells.Item(31,4).Value = "Giving shape to the public sector’s exemplary role in sustainable procurement"

$ws.Cells.Item(32,1).Value = "Z13_B01_P01"
$ws.Cells.Item(32,2).Value = "Z13_B01"
$ws.Cells.Item(32,3).Value = "Treibhausgase reduzieren"
$ws.Cells.Item(32,4).Value = "Reducing greenhouse gases"

$ws.Cells.Item(33,1).Value = "Z13_B01_P02"
$ws.Cells.Item(33,2).Value = "Z13_B01"
$ws.Cells.Item(33,3).Value = "Beitrag zur internationalen Klimafinanzierung leisten"
$ws.Cells.Item(33,4).Value = "Germany's contribution to international climate finance"

$ws.Cells.Item(34,1).Value = "Z14_B01_P01"
$ws.Cells.Item(34,2).Value = "Z14_B01"
$ws.Cells.Item(34,3).Value = "Meere und Meeresressourcen schützen und nachhaltig nutzen"
$ws.Cells.Item(34,4).Value = "Protecting and sustainably using oceans and marine resources"

$ws.Cells.Item(35,1).Value = "Z15_B01_P01"
$ws.Cells.Item(35,2).Value = "Z15_B01"
$ws.Cells.Item(35,3).Value = "Arten erhalten – Lebensräume schützen"
$ws.Cells.Item(35,4).Value = "Conserving species – Protecting habitats"

$ws.Cells.Item(36,1).Value = "Z15_B02_P01"
$ws.Cells.Item(36,2).Value = "Z15_B02"
$ws.Cells.Item(36,3).Value = "Ökosysteme schützen, Ökosystemleistungen erhalten und Lebensräume bewahren"
$ws.Cells.Item(36,4).Value = "Protecting ecosystems, conserving ecosystem services and preserving habitats"

$ws.Cells.Item(37,1).Value = "Z15_B02_P02"
$ws.Cells.Item(37,2).Value = "Z15_B02"
$ws.Cells.Item(37,3).Value = "Weltweit Entwaldung vermeiden und Böden schützen"
$ws.Cells.Item(37,4).Value = "Preventing deforestation and protecting soils world-wide"

$ws.Cells.Item(38,1).Value = "Z16_B01_P01"
$ws.Cells.Item(38,2).Value = "Z16_B01"
$ws.Cells.Item(38,3).Value = "Persönliche Sicherheit weiter erhöhen"
$ws.Cells.Item(38,4).Value = "Further increasing personal security"

$ws.Cells.Item(39,1).Value = "Z16_B02_P01"
$ws.Cells.Item(39,2).Value = "Z16_B02"
$ws.Cells.Item(39,3).Value = "Praktische Maßnahmen zur Bekämpfung der Proliferation, insbesondere von Kleinwaffen, ergreifen"
$ws.Cells.Item(39,4).Value = "Taking practical action to combat proliferation, especially of small arms"

$ws.Cells.Item(40,1).Value = "Z16_B03_P01"
$ws.Cells.Item(40,2).Value = "Z16_B03"
$ws.Cells.Item(40,3).Value = "Korruption bekämpfen"
$ws.Cells.Item(40,4).Value = "Combating corruption"

$ws.Cells.Item(41,1).Value = "Z17_B01_P01"
$ws.Cells.Item(41,2).Value = "Z17_B01"
$ws.Cells.Item(41,3).Value = "Nachhaltige Entwicklung unterstützen"
$ws.Cells.Item(41,4).Value = "Supporting sustainable development"

$ws.Cells.Item(42,1).Value = "Z17_B02_P01"
$ws.Cells.Item(42,2).Value = "Z17_B02"
$ws.Cells.Item(42,3).Value = "Wissen international vermitteln"
$ws.Cells.Item(42,4).Value = "Sharing knowledge internationally"

$ws.Cells.Item(43,1).Value = "Z17_B03_P01"
$ws.Cells.Item(43,2).Value = "Z17_B03"
$ws.Cells.Item(43,3).Value = "Handelschancen der Entwicklungsländer verbessern"
$ws.Cells.Item(43,4).Value = "Improving trade opportunities for developing countries"

